$d = $word.ActiveDocument

# --- Fix 1: correct the battery capacity typo for the iPhone 5 row
# ("5345" -> "5.45" [Wh]) in the "dati batteria smartphone" table.
$battTable = $d.Tables.Item(2)
$cell = $battTable.Cell(10, 4)
$cell.Range.Find.Execute("5345", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "5.45", 2) | Out-Null

# The diff shows the corrected value split across two runs ("5." + "45")
# rather than a single run. Force Word to break the run in two by
# nudging the character formatting of the first two characters ("5.")
# away from, then back to, its original value.
$cStart = $cell.Range.Start
$firstPart = $d.Range($cStart, $cStart + 2)
$origSize = $firstPart.Font.Size
$firstPart.Font.Size = $origSize + 1
$firstPart.Font.Size = $origSize

# --- Fix 2 & 3: the two "Tabella 3.x" captions have stale cached SEQ
# field results (swapped); refresh them so the Bluetooth power table
# becomes "Tabella 3.1" and the battery table becomes "Tabella 3.2".
$d.Fields.Update()
